# Daily price-data update: a new day's row is inserted at the top of the
# data table (row 2), pushing all existing rows down by one. The new row
# carries the latest date together with the same metric values that the
# rest of the (currently flat) series already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new day by inserting a row above the current row 2
# (the first data row, right below the header).
$ws.Rows.Item(2).Insert()

# Column A holds dates that are stored as plain text (e.g. "2026-02-19"),
# not real Excel date serials. Force the new cell to Text format before
# assigning the value so it is kept as a literal string instead of being
# auto-converted into a date number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-20"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
